# Append three new ticker rows to the bottom of the data range in Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newValues = @("IMX-USD", "TAO-USD", "GRT-USD")

# Find the last used row in column A and append after it.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

foreach ($val in $newValues) {
    $lastRow = $lastRow + 1
    $ws.Cells.Item($lastRow, 1).Value = $val
}
